# Auto-generated edit script applying the diff to Sheets/Typhon_Profits.xlsx
# (workbook contains sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1561.8405
$ws.Range("I15").Value = 1561.8405
$ws.Range("K15").Value = 4685.5215
$ws.Range("M15").Value = -4516.5215

$ws.Range("H43").Value = 1136.4615
$ws.Range("I43").Value = 1200
$ws.Range("J43").Value = 1131.1666
$ws.Range("K43").Value = 1200
$ws.Range("L43").Value = 1131.1666
$ws.Range("M43").Value = -1131
$ws.Range("N43").Value = -1269.1666

$ws.Range("H53").Value = 2547.3845
$ws.Range("I53").Value = 272.5
$ws.Range("J53").Value = 3558.4443
$ws.Range("K53").Value = 272.5
$ws.Range("L53").Value = 3558.4443
$ws.Range("M53").Value = 364.5
$ws.Range("N53").Value = -4832.4443

$ws.Range("H76").Value = 3475376.2
$ws.Range("I76").Value = 3364.2
$ws.Range("J76").Value = 55555556
$ws.Range("K76").Value = 3364.2
$ws.Range("L76").Value = 55555556
$ws.Range("M76").Value = -3049.2
$ws.Range("N76").Value = -55556186

$ws.Range("H79").Value = 3475376.2
$ws.Range("I79").Value = 3364.2
$ws.Range("J79").Value = 55555556
$ws.Range("K79").Value = 3364.2
$ws.Range("L79").Value = 55555556
$ws.Range("M79").Value = -2272.2
$ws.Range("N79").Value = -55557740

$ws.Range("H98").Value = 662.6
$ws.Range("I98").Value = 427.1875
$ws.Range("J98").Value = 1081.1111
$ws.Range("K98").Value = 427.1875
$ws.Range("L98").Value = 1081.1111
$ws.Range("M98").Value = 1070.8125
$ws.Range("N98").Value = -4077.1111

$ws.Range("H106").Value = 8549730
$ws.Range("I106").Value = 17545668
$ws.Range("J106").Value = 3590.4
$ws.Range("K106").Value = 17545668
$ws.Range("L106").Value = 3590.4
$ws.Range("M106").Value = -17545037
$ws.Range("N106").Value = -4852.4

$ws.Range("H112").Value = 1005.25
$ws.Range("J112").Value = 1094.4482
$ws.Range("L112").Value = 3283.3446
$ws.Range("N112").Value = -5499.3446

$ws.Range("H113").Value = 43481452
$ws.Range("I113").Value = 55558016
$ws.Range("J113").Value = 5820
$ws.Range("K113").Value = 55558016
$ws.Range("L113").Value = 5820
$ws.Range("M113").Value = -55554762
$ws.Range("N113").Value = -12328

$ws.Range("H116").Value = 9251.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 9251.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 9251.5
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -16135.5

$ws.Range("H122").Value = 662.6
$ws.Range("I122").Value = 427.1875
$ws.Range("J122").Value = 1081.1111
$ws.Range("K122").Value = 1281.5625
$ws.Range("L122").Value = 3243.3333
$ws.Range("M122").Value = 1168.4375
$ws.Range("N122").Value = -8143.3333

$ws.Range("H125").Value = 712.6667
$ws.Range("I125").Value = 819
$ws.Range("J125").Value = 500
$ws.Range("K125").Value = 7371
$ws.Range("L125").Value = 4500
$ws.Range("M125").Value = -4911
$ws.Range("N125").Value = -9420

$ws.Range("H129").Value = 176479.4
$ws.Range("J129").Value = 186265.36
$ws.Range("L129").Value = 558796.08
$ws.Range("N129").Value = -568796.08

$ws.Range("H135").Value = 17246966
$ws.Range("I135").Value = 692.2381
$ws.Range("J135").Value = 62518430
$ws.Range("K135").Value = 6230.142900000001
$ws.Range("L135").Value = 562665870
$ws.Range("M135").Value = -3695.142900000001
$ws.Range("N135").Value = -562670940

$ws.Range("H137").Value = 1501.9445
$ws.Range("I137").Value = 1516.909
$ws.Range("J137").Value = 1478.4286
$ws.Range("K137").Value = 4550.727000000001
$ws.Range("L137").Value = 4435.2858
$ws.Range("M137").Value = -2000.727000000001
$ws.Range("N137").Value = -9535.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11509.48
$ws.Range("I32").Value = 10874.111
$ws.Range("J32").Value = 13143.286
$ws.Range("K32").Value = 10874.111
$ws.Range("L32").Value = 13143.286
$ws.Range("M32").Value = -10587.111
$ws.Range("N32").Value = -13717.286

$ws.Range("H45").Value = 3529.2856
$ws.Range("I45").Value = 2956.8572
$ws.Range("K45").Value = 2956.8572
$ws.Range("M45").Value = -2579.8572

$ws.Range("H61").Value = 1373.375
$ws.Range("I61").Value = 1322.2174
$ws.Range("J61").Value = 2550
$ws.Range("K61").Value = 1322.2174
$ws.Range("L61").Value = 2550
$ws.Range("M61").Value = -1110.2174
$ws.Range("N61").Value = -2974

$ws.Range("H136").Value = 1373.375
$ws.Range("I136").Value = 1322.2174
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 3966.6522
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = -1416.6522
$ws.Range("N136").Value = -12750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 16000
$ws.Range("J56").Value = 16000
$ws.Range("L56").Value = 16000
$ws.Range("N56").Value = -17478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()

$ws.Range("H105").Value = 811.9583
$ws.Range("I105").Value = 743.619
$ws.Range("K105").Value = 743.619
$ws.Range("M105").Value = 1003.381

$ws.Range("H122").Value = 1064
$ws.Range("I122").Value = 926.8
$ws.Range("K122").Value = 2780.4
$ws.Range("M122").Value = -330.3999999999996

$ws.Range("H125").Value = 10300
$ws.Range("I125").Value = 5500
$ws.Range("J125").Value = 19900
$ws.Range("K125").Value = 5500
$ws.Range("L125").Value = 19900
$ws.Range("M125").Value = -3040
$ws.Range("N125").Value = -24820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4522.0586
$ws.Range("I63").Value = 1803
$ws.Range("J63").Value = 5358.6924
$ws.Range("K63").Value = 5409
$ws.Range("L63").Value = 16076.0772
$ws.Range("M63").Value = -4660
$ws.Range("N63").Value = -17574.0772

$ws.Range("H64").Value = 4188
$ws.Range("I64").Value = 1812
$ws.Range("J64").Value = 4425.6
$ws.Range("K64").Value = 5436
$ws.Range("L64").Value = 13276.8
$ws.Range("M64").Value = -5166
$ws.Range("N64").Value = -13816.8

$ws.Range("H66").Value = 4522.0586
$ws.Range("I66").Value = 1803
$ws.Range("J66").Value = 5358.6924
$ws.Range("K66").Value = 16227
$ws.Range("L66").Value = 48228.2316
$ws.Range("M66").Value = -12483
$ws.Range("N66").Value = -55716.2316

$ws.Range("H67").Value = 4188
$ws.Range("I67").Value = 1812
$ws.Range("J67").Value = 4425.6
$ws.Range("K67").Value = 5436
$ws.Range("L67").Value = 13276.8
$ws.Range("M67").Value = -4500
$ws.Range("N67").Value = -15148.8

$ws.Range("H68").Value = 1790.5834
$ws.Range("J68").Value = 1931.8889
$ws.Range("L68").Value = 5795.6667
$ws.Range("N68").Value = -7417.6667

$ws.Range("H71").Value = 1790.5834
$ws.Range("J71").Value = 1931.8889
$ws.Range("L71").Value = 17387.0001
$ws.Range("N71").Value = -25499.0001

$ws.Range("H107").Value = 5177.1
$ws.Range("I107").Value = 5726.8335
$ws.Range("K107").Value = 17180.5005
$ws.Range("M107").Value = -15260.5005

$ws.Range("H122").Value = 733.58826
$ws.Range("I122").Value = 338.8
$ws.Range("K122").Value = 3049.2
$ws.Range("M122").Value = -599.2000000000003

$ws.Range("H131").Value = 119847.5
$ws.Range("I131").Value = 760
$ws.Range("J131").Value = 124258.15
$ws.Range("K131").Value = 2280
$ws.Range("L131").Value = 372774.45
$ws.Range("M131").Value = 2760
$ws.Range("N131").Value = -382854.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 28296.666
$ws.Range("J57").Value = 28296.666
$ws.Range("L57").Value = 28296.666
$ws.Range("N57").Value = -29936.666

$ws.Range("H102").Value = 26318908
$ws.Range("I102").Value = 31252886
$ws.Range("J102").Value = 4354.6665
$ws.Range("K102").Value = 31252886
$ws.Range("L102").Value = 4354.6665
$ws.Range("M102").Value = -31251264
$ws.Range("N102").Value = -7598.6665

$ws.Range("H126").Value = 5710.2964
$ws.Range("I126").Value = 4615.143
$ws.Range("K126").Value = 13845.429
$ws.Range("M126").Value = -11375.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4408.25
$ws.Range("I7").Value = 4457.143
$ws.Range("K7").Value = 4457.143
$ws.Range("M7").Value = -4345.143

$ws.Range("H40").Value = 3285.2
$ws.Range("I40").Value = 2748
$ws.Range("J40").Value = 4426.75
$ws.Range("K40").Value = 2748
$ws.Range("L40").Value = 4426.75
$ws.Range("M40").Value = -2612
$ws.Range("N40").Value = -4698.75

$ws.Range("H46").Value = 1022.58826
$ws.Range("I46").Value = 996
$ws.Range("J46").Value = 1900
$ws.Range("K46").Value = 996
$ws.Range("L46").Value = 1900
$ws.Range("M46").Value = -808
$ws.Range("N46").Value = -2276

$ws.Range("H126").Value = 4408.25
$ws.Range("I126").Value = 4457.143
$ws.Range("K126").Value = 13371.429
$ws.Range("M126").Value = -10901.429

$ws.Range("H136").Value = 2519.6
$ws.Range("I136").Value = 2519.6
$ws.Range("K136").Value = 7558.799999999999
$ws.Range("M136").Value = -5008.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 8936
$ws.Range("I64").Value = 8936
$ws.Range("K64").Value = 8936
$ws.Range("M64").Value = -8688

$ws.Range("H67").Value = 8936
$ws.Range("I67").Value = 8936
$ws.Range("K67").Value = 8936
$ws.Range("M67").Value = -8078

$ws.Range("H126").Value = 1246.0667
$ws.Range("I126").Value = 1246.0667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3738.2001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1268.2001
$ws.Range("N126").ClearContents()
